$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Paragraph 2 of the "Content Placeholder 2" text frame currently reads:
#   "Cong viec du dinh lam hom truoc:"
# and must become:
#   "Cong viec da du dinh lam:"
#
# Apply the word-level substitutions (and the trailing deletion) from the
# right-hand side of the paragraph towards the left so that earlier
# (left-most) character offsets, computed against the original text,
# remain valid while later (right-most) edits are performed first.

$para = $tr.Paragraphs(2, 1)

# Remove the trailing " truoc" (space + word) - characters 26..31 of the
# original paragraph text - leaving the final ":" run untouched.
$para.Characters(26, 6).Text = ""

# "hom" (chars 23..25) becomes "lam"
$para.Characters(23, 3).Text = "làm"

# "lam" (chars 19..21) becomes "dinh"
$para.Characters(19, 3).Text = "định"

# "dinh" (chars 14..17) becomes "du"
$para.Characters(14, 4).Text = "dự"

# "du" (chars 11..12) becomes "da"
$para.Characters(11, 2).Text = "đã"
